$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.703.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.164.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.31"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +11.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.423"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.90%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.708.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.52%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.760.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.158.43"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.47"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +9.87%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.77"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.47"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.518"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.95%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +14.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0880"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.82%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.96"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.34"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.38"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +13.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.54"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.650.69"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.89%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.71"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.708"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0279"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.95%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +14.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.23"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.984"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.30"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.760"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.61%  "
